# "Generate Report for Handoff"
# Updates the localization-status report: flips the in-progress status to
# "Ready for handoff" and refreshes the handoff timestamps, then widens the
# datetime columns that now need to fit the longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Latest handoff / HO xliff generate datetimes ---
$wsOverview.Range("G2").Value = "2016-08-21 23:06:01"
$wsDeDe.Range("H2").Value = "2016-08-21 23:06:01"
$wsZhCn.Range("H2").Value = "2016-08-21 23:05:56"

# --- Widen the datetime columns to fit the new status/text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
